$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input Values")

# Selenium Grid config now accepts the "ANY" browser (was hard-coded to
# Chrome) so the Firefox-driven run can execute against this test data.
$ws.Range("B1").Value = "ANY"

# Keep the active selection on the edited cell, matching the authored state.
$ws.Activate()
$ws.Range("B1").Select()
